# Weekly update: insert a new price-report row for
# "Comercializadora del Agro de Limarí - Poroto granado" ahead of the
# existing row 51, shifting the rest of the series down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 51; everything from 51..124 shifts to 52..125
# and the sheet's used range grows to A1:R125.
$ws.Rows("51:51").Insert()

# Populate the newly inserted row with the latest week's data.
$ws.Range("A51").Value = 2
$ws.Range("B51").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C51").Value = "Coquimbo"
$ws.Range("D51").Value = 45036
$ws.Range("E51").Value = 4
$ws.Range("F51").Value = 100112030
$ws.Range("G51").Value = "Poroto granado"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 400
$ws.Range("K51").Value = 27000
$ws.Range("L51").Value = 28000
$ws.Range("M51").Value = 27500
$ws.Range("N51").Value = "$/malla 25 kilos"
$ws.Range("O51").Value = "Provincia de Limarí"
$ws.Range("P51").Value = 1100
$ws.Range("Q51").Value = 25
$ws.Range("R51").Value = "Hortaliza"
